$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "2024-07-30 Tuesday" "2024-07-31 Wednesday"

Replace-Text "40÷6=6, 4" "28÷9=3, 1"
Replace-Text "64÷3=21, 1" "39÷5=7, 4"
Replace-Text "52÷2=26, 0" "48÷3=16, 0"
Replace-Text "67÷3=22, 1" "98÷5=19, 3"
Replace-Text "84÷2=42, 0" "47÷5=9, 2"

Replace-Text "47÷2=23, 1" "98÷8=12, 2"
Replace-Text "75÷8=9, 3" "13÷6=2, 1"
Replace-Text "90÷7=12, 6" "16÷7=2, 2"
Replace-Text "86÷6=14, 2" "18÷4=4, 2"
Replace-Text "25÷6=4, 1" "43÷8=5, 3"

Replace-Text "12÷5=2, 2" "77÷6=12, 5"
Replace-Text "74÷7=10, 4" "25÷5=5, 0"
Replace-Text "64÷4=16, 0" "93÷9=10, 3"
Replace-Text "84÷4=21, 0" "92÷6=15, 2"
Replace-Text "43÷2=21, 1" "21÷4=5, 1"

Replace-Text "31÷9=3, 4" "86÷4=21, 2"
Replace-Text "13÷3=4, 1" "35÷2=17, 1"
Replace-Text "91÷5=18, 1" "58÷2=29, 0"
Replace-Text "65÷2=32, 1" "53÷8=6, 5"
Replace-Text "12÷3=4, 0" "47÷5=9, 2"

Replace-Text "20÷2=10, 0" "46÷6=7, 4"
Replace-Text "69÷2=34, 1" "56÷7=8, 0"
Replace-Text "68÷5=13, 3" "28÷7=4, 0"
Replace-Text "33÷7=4, 5" "74÷6=12, 2"
Replace-Text "69÷3=23, 0" "30÷8=3, 6"
